# Adds four new species-observation rows (3-6) to the "Artfynd" sheet,
# mirroring the rows appended upstream. Numbers are written as numbers,
# booleans as booleans, and text that Excel would otherwise auto-convert
# (pure digits, ISO dates) is entered with a leading apostrophe and then
# has ClearFormats() applied so it stays plain text without leaving a
# quote-prefix style behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 3
$ws.Cells.Item(3,1).Value = 112313702
$ws.Cells.Item(3,2).Value = 90800
$ws.Cells.Item(3,3).Value = 'Ovaliderad'
$ws.Cells.Item(3,4).Value = 'LC'
$ws.Cells.Item(3,5).Value = 4364
$ws.Cells.Item(3,6).Value = 'Dropptaggsvamp'
$ws.Cells.Item(3,7).Value = 'Hydnellum ferrugineum'
$ws.Cells.Item(3,8).Value = '(Fr.:Fr.) P. Karst.'
# (row 3, col 9 in the source is an explicit empty-text cell; COM normalises any "" write back to a truly blank cell, so it is intentionally left unset here)
$ws.Cells.Item(3,16).Value = 'Trollberget-Kullaheden, Vb'
$ws.Cells.Item(3,17).Value = 752827
$ws.Cells.Item(3,18).Value = 7093488
$ws.Cells.Item(3,19).Value = 10
$ws.Cells.Item(3,20).Value = 'Västerbotten'
$ws.Cells.Item(3,21).Value = 'Umeå'
$ws.Cells.Item(3,22).Value = 'Västerbotten'
$ws.Cells.Item(3,23).Value = 'Umeå socken'
$c = $ws.Cells.Item(3,25)
$c.Value = "'2023-09-24"
$c.ClearFormats()
$c = $ws.Cells.Item(3,27)
$c.Value = "'2023-09-24"
$c.ClearFormats()
$ws.Cells.Item(3,30).Value = $false
$ws.Cells.Item(3,31).Value = $false
$ws.Cells.Item(3,33).Value = $false
# (row 3, col 46 in the source is an explicit empty-text cell; COM normalises any "" write back to a truly blank cell, so it is intentionally left unset here)
$ws.Cells.Item(3,49).Value = 'Roger Olofsson'
$ws.Cells.Item(3,50).Value = 'Roger Olofsson'
# (row 3, col 51 in the source is an explicit empty-text cell; COM normalises any "" write back to a truly blank cell, so it is intentionally left unset here)

# Row 4
$ws.Cells.Item(4,1).Value = 112313576
$ws.Cells.Item(4,2).Value = 56575
$ws.Cells.Item(4,3).Value = 'Ovaliderad'
$ws.Cells.Item(4,4).Value = 'NT'
$ws.Cells.Item(4,5).Value = 103021
$ws.Cells.Item(4,6).Value = 'Talltita'
$ws.Cells.Item(4,7).Value = 'Poecile montanus'
$ws.Cells.Item(4,8).Value = '(Conrad von Baldenstein, 1827)'
$c = $ws.Cells.Item(4,9)
$c.Value = "'1"
$c.ClearFormats()
$ws.Cells.Item(4,13).Value = 'lockläte, övriga läten'
$ws.Cells.Item(4,16).Value = 'Trollberget-Kullaheden, Vb'
$ws.Cells.Item(4,17).Value = 752714
$ws.Cells.Item(4,18).Value = 7093570
$ws.Cells.Item(4,19).Value = 10
$ws.Cells.Item(4,20).Value = 'Västerbotten'
$ws.Cells.Item(4,21).Value = 'Umeå'
$ws.Cells.Item(4,22).Value = 'Västerbotten'
$ws.Cells.Item(4,23).Value = 'Umeå socken'
$c = $ws.Cells.Item(4,25)
$c.Value = "'2023-09-24"
$c.ClearFormats()
$c = $ws.Cells.Item(4,27)
$c.Value = "'2023-09-24"
$c.ClearFormats()
$ws.Cells.Item(4,30).Value = $false
$ws.Cells.Item(4,31).Value = $false
$ws.Cells.Item(4,33).Value = $false
# (row 4, col 46 in the source is an explicit empty-text cell; COM normalises any "" write back to a truly blank cell, so it is intentionally left unset here)
$ws.Cells.Item(4,49).Value = 'Roger Olofsson'
$ws.Cells.Item(4,50).Value = 'Roger Olofsson'
# (row 4, col 51 in the source is an explicit empty-text cell; COM normalises any "" write back to a truly blank cell, so it is intentionally left unset here)

# Row 5
$ws.Cells.Item(5,1).Value = 112313655
$ws.Cells.Item(5,2).Value = 90823
$ws.Cells.Item(5,3).Value = 'Ovaliderad'
$ws.Cells.Item(5,4).Value = 'NT'
$ws.Cells.Item(5,5).Value = 5966
$ws.Cells.Item(5,6).Value = 'Motaggsvamp'
$ws.Cells.Item(5,7).Value = 'Sarcodon squamosus'
$ws.Cells.Item(5,8).Value = '(Schaeff.) Quél.'
$c = $ws.Cells.Item(5,9)
$c.Value = "'2"
$c.ClearFormats()
$ws.Cells.Item(5,10).Value = 'fruktkroppar'
$ws.Cells.Item(5,16).Value = 'Trollberget-Kullaheden, Vb'
$ws.Cells.Item(5,17).Value = 752820
$ws.Cells.Item(5,18).Value = 7093493
$ws.Cells.Item(5,19).Value = 10
$ws.Cells.Item(5,20).Value = 'Västerbotten'
$ws.Cells.Item(5,21).Value = 'Umeå'
$ws.Cells.Item(5,22).Value = 'Västerbotten'
$ws.Cells.Item(5,23).Value = 'Umeå socken'
$c = $ws.Cells.Item(5,25)
$c.Value = "'2023-09-24"
$c.ClearFormats()
$c = $ws.Cells.Item(5,27)
$c.Value = "'2023-09-24"
$c.ClearFormats()
$ws.Cells.Item(5,30).Value = $false
$ws.Cells.Item(5,31).Value = $false
$ws.Cells.Item(5,33).Value = $false
# (row 5, col 46 in the source is an explicit empty-text cell; COM normalises any "" write back to a truly blank cell, so it is intentionally left unset here)
$ws.Cells.Item(5,49).Value = 'Roger Olofsson'
$ws.Cells.Item(5,50).Value = 'Roger Olofsson'
# (row 5, col 51 in the source is an explicit empty-text cell; COM normalises any "" write back to a truly blank cell, so it is intentionally left unset here)

# Row 6
$ws.Cells.Item(6,1).Value = 112313590
$ws.Cells.Item(6,2).Value = 56446
$ws.Cells.Item(6,3).Value = 'Ovaliderad'
$ws.Cells.Item(6,4).Value = 'NT'
$ws.Cells.Item(6,5).Value = 100049
$ws.Cells.Item(6,6).Value = 'Spillkråka'
$ws.Cells.Item(6,7).Value = 'Dryocopus martius'
$ws.Cells.Item(6,8).Value = '(Linnaeus, 1758)'
# (row 6, col 9 in the source is an explicit empty-text cell; COM normalises any "" write back to a truly blank cell, so it is intentionally left unset here)
$ws.Cells.Item(6,13).Value = 'äldre spår'
$ws.Cells.Item(6,16).Value = 'Trollberget-Kullaheden, Vb'
$ws.Cells.Item(6,17).Value = 752543
$ws.Cells.Item(6,18).Value = 7093684
$ws.Cells.Item(6,19).Value = 10
$ws.Cells.Item(6,20).Value = 'Västerbotten'
$ws.Cells.Item(6,21).Value = 'Umeå'
$ws.Cells.Item(6,22).Value = 'Västerbotten'
$ws.Cells.Item(6,23).Value = 'Umeå socken'
$c = $ws.Cells.Item(6,25)
$c.Value = "'2023-09-24"
$c.ClearFormats()
$c = $ws.Cells.Item(6,27)
$c.Value = "'2023-09-24"
$c.ClearFormats()
$ws.Cells.Item(6,30).Value = $false
$ws.Cells.Item(6,31).Value = $false
$ws.Cells.Item(6,33).Value = $false
# (row 6, col 46 in the source is an explicit empty-text cell; COM normalises any "" write back to a truly blank cell, so it is intentionally left unset here)
$ws.Cells.Item(6,49).Value = 'Roger Olofsson'
$ws.Cells.Item(6,50).Value = 'Roger Olofsson'
# (row 6, col 51 in the source is an explicit empty-text cell; COM normalises any "" write back to a truly blank cell, so it is intentionally left unset here)
